$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.599.31"
$ws.Range("E2").Value = "  +1.47%  "

$ws.Range("D3").Value = "2.247.51"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("D4").Value = "'1.02"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.31%  "

$ws.Range("D5").Value = "'309.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("D6").Value = "'94.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").Value = "'0.575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.37%  "

$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").Value = "'0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").Value = "'34.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").Value = "'0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").Value = "'7.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.329.42"
$ws.Range("E14").Value = "  +4.22%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.840"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.31%  "

$ws.Range("D16").Value = "'13.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").Value = "44.307.94"

$ws.Range("D18").Value = "0.0₃0961"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.10%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "'65.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "

$ws.Range("D22").Value = "'238.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").Value = "'2.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.44%  "

$ws.Range("D24").Value = "'2.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.56%  "

$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").Value = "'2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.93%  "

$ws.Range("D27").Value = "'9.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("D28").Value = "'37.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.04%  "

$ws.Range("D29").Value = "'6.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("D30").Value = "'20.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.47%  "

$ws.Range("D31").Value = "'152.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "

$ws.Range("D32").Value = "'0.0806"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("D34").Value = "'3.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.73%  "

$ws.Range("D35").Value = "'0.110"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.04%  "

$ws.Range("D36").Value = "'0.120"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.61%  "

$ws.Range("D37").Value = "'1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.16%  "

$ws.Range("D38").Value = "'3.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.77%  "

$ws.Range("D39").Value = "'3.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").Value = "'14.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.75%  "

$ws.Range("D41").Value = "'0.0302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("E42").Value = "  +0.31%  "

$ws.Range("D43").Value = "1.747.69"
$ws.Range("E43").Value = "  +0.84%  "

$ws.Range("D44").Value = "'0.194"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.90%  "

$ws.Range("D45").Value = "'80.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.70%  "

$ws.Range("D46").Value = "'99.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "

$ws.Range("D47").Value = "'70.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.51%  "

$ws.Range("D48").Value = "'4.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "'55.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.18%  "

$ws.Range("D50").Value = "'8.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.87%  "

$ws.Range("D51").Value = "'1.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.40%  "
